# Update "paises" (countries) COVID tracking sheet:
#  - South Africa (Sudafrica) overtakes United Kingdom (Reino Unido) in the ranking
#  - Costa Rica overtakes Republica de Macedonia / Senegal / Consejo Danes para los
#    Refugiados in the ranking (each of those three countries is pushed down one row)
#  - Refresh numeric stats for several rows
#  - Bump the "last updated" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos (stats refresh only, still #1) ---
$ws.Range("B4").Value = 3518037
$ws.Range("C4").Value = 38554
$ws.Range("D4").Value = 1568330
$ws.Range("E4").Value = 1810909
$ws.Range("G4").Value = 551
$ws.Range("H4").Value = 138798

# --- Rows 12-13: Sudafrica overtakes Reino Unido ---
$ws.Range("A12").Value = "Sudafrica"
$ws.Range("B12").Value = 298292
$ws.Range("C12").Value = 10496
$ws.Range("D12").Value = 146279
$ws.Range("E12").Value = 147667
$ws.Range("G12").Value = 174
$ws.Range("H12").Value = 4346

$ws.Range("A13").Value = "Reino Unido"
$ws.Range("B13").Value = 291373
$ws.Range("C13").Value = 398
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("G13").Value = 138
$ws.Range("H13").Value = 44968

# --- Row 19: Alemania (stats refresh only) ---
$ws.Range("B19").Value = 200704
$ws.Range("C19").Value = 268
$ws.Range("D19").Value = 185500
$ws.Range("E19").Value = 6063
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 9141

# --- Rows 80-83: Costa Rica overtakes Republica de Macedonia, Senegal and
#     Consejo Danes para los Refugiados (each pushed down one row) ---
$ws.Range("A80").Value = "Costa Rica"
$ws.Range("B80").Value = 8482
$ws.Range("C80").Value = 446
$ws.Range("D80").Value = 2441
$ws.Range("E80").Value = 6005
$ws.Range("G80").Value = 5
$ws.Range("H80").Value = 36

$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("B81").Value = 8332
$ws.Range("C81").Value = 135
$ws.Range("D81").Value = 4468
$ws.Range("E81").Value = 3475
$ws.Range("G81").Value = 4
$ws.Range("H81").Value = 389

$ws.Range("A82").Value = "Senegal"
$ws.Range("B82").Value = 8243
$ws.Range("C82").Value = 45
$ws.Range("D82").Value = 5580
$ws.Range("E82").Value = 2513
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 150

$ws.Range("A83").Value = "Consejo Danes para los Refugiados"
$ws.Range("B83").Value = 8135
$ws.Range("C83").Value = 60
$ws.Range("D83").Value = 3948
$ws.Range("E83").Value = 3997
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 190

# --- Row 94: Mauritania (stats refresh only) ---
$ws.Range("B94").Value = 5518
$ws.Range("C94").Value = 72
$ws.Range("D94").Value = 2664
$ws.Range("E94").Value = 2707

# --- Row 104: Somalia (stats refresh only) ---
$ws.Range("B104").Value = 3076
$ws.Range("C104").Value = 4
$ws.Range("D104").Value = 1380
$ws.Range("E104").Value = 1603

# --- Row 148: Principado de Andorra (stats refresh only) ---
$ws.Range("B148").Value = 861
$ws.Range("C148").Value = 3
$ws.Range("E148").Value = 6

# --- Row 149: Surinam (stats refresh only) ---
$ws.Range("B149").Value = 801
$ws.Range("C149").Value = 21
$ws.Range("D149").Value = 543
$ws.Range("E149").Value = 240

# --- Row 183: Aruba (stats refresh only) ---
$ws.Range("B183").Value = 106
$ws.Range("C183").Value = 1
$ws.Range("E183").Value = 4

# --- Bump "last updated" timestamp (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 21:28"
